# Applies the "Added calcium and b12 to Zambia and Uganda and updated output"
# refresh to the sessionInfo/Info sheets of the workbook:
#   - Info sheet: Start_time / End_time timestamps updated (re-run of the script)
#   - sessionInfo sheet: a handful of R package version numbers were bumped
#     ("here", "magrittr", "rprojroot"), and the "backports" package dropped
#     out of the "Loaded only" package list (so "boot", which used to be the
#     entry right after it, shifts up one row and the now-superfluous last
#     row of that list is cleared).

$wb = $excel.ActiveWorkbook

# ---- Info sheet: refresh run timestamps ----
$wsInfo = $wb.Worksheets.Item("Info")
$wsInfo.Range("B26").Value = "Thu Nov 19 15:24:00 2020"   # Start_time
$wsInfo.Range("B27").Value = "Thu Nov 19 15:24:05 2020"   # End_time

# ---- sessionInfo sheet: refresh package version numbers ----
$wsSession = $wb.Worksheets.Item("sessionInfo")

$wsSession.Range("G2").Value  = "1.0.0"    # here
$wsSession.Range("J3").Value  = "2.0.1"    # magrittr
$wsSession.Range("J10").Value = "2.0.2"    # rprojroot

# "backports" (previously on row 15 of the Loaded-only list) is no longer
# loaded; "boot" (previously row 16) moves up to row 15, and row 16's
# entry for that column pair is cleared out entirely.
$wsSession.Range("I15").Value = "boot"
$wsSession.Range("J15").Value = "1.3-25"
$wsSession.Range("I16:J16").ClearContents()
